$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the formulas (and cached values) out of C2:C5, keep the existing style.
$ws.Range("C2:C5").ClearContents()

# Append new rows of data.
$ws.Range("A6").Value = "Opinion Trading"
$ws.Range("B6").Value = "html/op_trading.html"
$ws.Range("A7").Value = "Madhuri Gupta"
$ws.Range("A8").Value = "Mohini Mohan Dutta"
$ws.Range("A9").Value = "Karnataka Soaps and Detergents Limited"
$ws.Range("A10").Value = "IndiGo"
$ws.Range("A11").Value = "Vaibhav Taneja"
$ws.Range("A12").Value = "Emirates Draw"

# Resize column A to fit the new, wider content (best-fit for the longest
# new entry, "Karnataka Soaps and Detergents Limited").
$ws.Columns.Item(1).ColumnWidth = 33.5

# Move the selection to reflect where the user left off editing.
$ws.Range("A13").Select()
